# The deck currently uses the "Integral" design theme (ppt/theme/theme1.xml,
# linked from the one slide master). The authored change swaps in the
# default Office color palette for that theme (the deck's XML also swaps
# the raw theme1.xml / theme2.xml parts so the notes master - which was
# already on "Office Theme" - ends up carrying the old "Integral" colors,
# but the PowerPoint object model only exposes the live/active theme that
# backs the slide master, so we reproduce the user-visible effect: the
# slide master's theme color scheme becomes the standard Office palette).

$p = $ppt.ActivePresentation

# Any slide's ThemeColorScheme reaches into the single shared theme used
# by the slide master/layouts, so slide 1 is as good an anchor as any.
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-ThemeColor($scheme, $idx, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $scheme.Colors($idx).RGB = $r + ($g * 256) + ($b * 65536)
}

# Index order matches the standard 12-slot theme colour scheme:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
Set-ThemeColor $tcs 1  "000000"
Set-ThemeColor $tcs 2  "FFFFFF"
Set-ThemeColor $tcs 3  "44546A"
Set-ThemeColor $tcs 4  "E7E6E6"
Set-ThemeColor $tcs 5  "5B9BD5"
Set-ThemeColor $tcs 6  "ED7D31"
Set-ThemeColor $tcs 7  "A5A5A5"
Set-ThemeColor $tcs 8  "FFC000"
Set-ThemeColor $tcs 9  "4472C4"
Set-ThemeColor $tcs 10 "70AD47"
Set-ThemeColor $tcs 11 "0563C1"
Set-ThemeColor $tcs 12 "954F72"
